$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and label (E) cells to match the refreshed symbol list scrape
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "271.40"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.68"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.309"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06293"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.557"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.568"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.374"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8244"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01378"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1584"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08337"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03398"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03194"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.060"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09242"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001667"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04689"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006266"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005971"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001063"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001490"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.763"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.367"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3322"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1254"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002706"
$ws.Range("E28").Value = "27UpBotsUBXT"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04713"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007060"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1170"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003632"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01167"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006019"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009830"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000744"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7763"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002327"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00001290"
$ws.Range("E50").Value = "49CryptobidCoinCBCWorstin24h"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01231"
